$wb = $excel.ActiveWorkbook

# --- zh-cn sheet (Row 3 = b105cc69 entry) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-28 16:16:13"

# --- de-de sheet (Row 3 = b105cc69 entry) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H3").Value = "2016-08-28 16:16:17"

# --- Overview sheet (Row 3 = b105cc69.md entry) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 16:16:17"
